# writeToExcel Hooks - update testCitizen sheet contents
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testCitizen")

# Column B must be (re)written before column A so the shared-string table
# is rebuilt in the same order as the target workbook (urbsNN4 block first,
# then the mutated ulaisNNNN block).
$ws.Range("B1").Value = "urbs134"
$ws.Range("B2").Value = "urbs144"
$ws.Range("B3").Value = "urbs154"
$ws.Range("B4").Value = "urbs164"
$ws.Range("B5").Value = "urbs174"
$ws.Range("B6").Value = "urbs184"
$ws.Range("B7").Value = "urbs194"
$ws.Range("B8").Value = "urbs204"

$ws.Range("A1").Value = "ul1ais1551"
$ws.Range("A2").Value = "ula2is11gg1"
$ws.Range("A3").Value = "ulai3s11ff1"
$ws.Range("A4").Value = "ulaise114ss1"
$ws.Range("A5").Value = "ulaise11xx1"
$ws.Range("A6").Value = "ulaeis1vv1"
$ws.Range("A7").Value = "ulaifs115bb1"
$ws.Range("A8").Value = "ulaisf11ff1"

# Column A got a touch wider
$ws.Columns.Item(1).ColumnWidth = 19.5

# Selection moved from D8 to A10
$ws.Activate()
[void]$ws.Range("A10").Select()
